$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("poc_config")

# "remove tokyo from all kpis": every data row (2-10) carried an
# address_city/Tokyo pair in columns J and K. Clear those cells (content
# only, keep formatting/style) for all rows so the shared strings
# "address_city" and "Tokyo" are dropped entirely from the workbook.
$ws.Range("J2:K10").ClearContents()

# Keep the selection/view consistent with where the edit happened.
$ws.Range("I32").Select() | Out-Null
